$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.7310083333333334
$ws.Range("H2").Value = 2.193025
$ws.Range("I2").Value = 0.01673731480740535
$ws.Range("J2").Value = 0.01673731480740535
$ws.Range("M2").Value = 0.3620403333333334
$ws.Range("N2").Value = 1.086121
$ws.Range("O2").Value = 0.138293228945037
$ws.Range("P2").Value = 0.138293228945037
$ws.Range("Q2").Value = 0.2646545006694445
$ws.Range("R2").Value = 2.381890506025
$ws.Range("S2").Value = 0.002314657308585666
$ws.Range("T2").Value = 0.002314657308585666
$ws.Range("G3").Value = 0.7310083333333334
$ws.Range("H3").Value = 2.193025
$ws.Range("I3").Value = 0.01673731480740535
$ws.Range("J3").Value = 0.01673731480740535
$ws.Range("O3").Value = 0.530310855165568
$ws.Range("P3").Value = 0.530310855165568
$ws.Range("Q3").Value = 1.014866422919444
$ws.Range("R3").Value = 9.133797806275
$ws.Range("S3").Value = 0.008875979728690454
$ws.Range("T3").Value = 0.008875979728690454
$ws.Range("G4").Value = 0.7310083333333334
$ws.Range("H4").Value = 2.193025
$ws.Range("I4").Value = 0.01673731480740535
$ws.Range("J4").Value = 0.01673731480740535
$ws.Range("M4").Value = 0.8675673333333332
$ws.Range("N4").Value = 2.602702
$ws.Range("O4").Value = 0.331395915889395
$ws.Range("P4").Value = 0.331395915889395
$ws.Range("Q4").Value = 0.6341989503944444
$ws.Range("R4").Value = 5.70779055355
$ws.Range("S4").Value = 0.005546677770129229
$ws.Range("T4").Value = 0.005546677770129229
$ws.Range("I5").Value = 0.8536212576586365
$ws.Range("J5").Value = 0.8536212576586365
$ws.Range("M5").Value = 0.3620403333333334
$ws.Range("N5").Value = 1.086121
$ws.Range("O5").Value = 0.138293228945037
$ws.Range("P5").Value = 0.138293228945037
$ws.Range("Q5").Value = 13.49766735620667
$ws.Range("R5").Value = 121.47900620586
$ws.Range("S5").Value = 0.1180500400177362
$ws.Range("T5").Value = 0.1180500400177362
$ws.Range("I6").Value = 0.8536212576586365
$ws.Range("J6").Value = 0.8536212576586365
$ws.Range("O6").Value = 0.530310855165568
$ws.Range("P6").Value = 0.530310855165568
$ws.Range("Q6").Value = 51.75929127560666
$ws.Range("S6").Value = 0.4526846191364592
$ws.Range("T6").Value = 0.4526846191364592
$ws.Range("I7").Value = 0.8536212576586365
$ws.Range("J7").Value = 0.8536212576586365
$ws.Range("M7").Value = 0.8675673333333332
$ws.Range("N7").Value = 2.602702
$ws.Range("O7").Value = 0.331395915889395
$ws.Range("P7").Value = 0.331395915889395
$ws.Range("Q7").Value = 32.34483618614666
$ws.Range("R7").Value = 291.1035256753199
$ws.Range("S7").Value = 0.2828865985044411
$ws.Range("T7").Value = 0.2828865985044411
$ws.Range("G8").Value = 5.662136666666666
$ws.Range("H8").Value = 16.98641
$ws.Range("I8").Value = 0.129641427533958
$ws.Range("J8").Value = 0.129641427533958
$ws.Range("M8").Value = 0.3620403333333334
$ws.Range("N8").Value = 1.086121
$ws.Range("O8").Value = 0.138293228945037
$ws.Range("P8").Value = 0.138293228945037
$ws.Range("Q8").Value = 2.049921846178889
$ws.Range("R8").Value = 18.44929661561
$ws.Range("S8").Value = 0.01792853161871508
$ws.Range("T8").Value = 0.01792853161871508
$ws.Range("G9").Value = 5.662136666666666
$ws.Range("H9").Value = 16.98641
$ws.Range("I9").Value = 0.129641427533958
$ws.Range("J9").Value = 0.129641427533958
$ws.Range("O9").Value = 0.530310855165568
$ws.Range("P9").Value = 0.530310855165568
$ws.Range("Q9").Value = 7.860802843078888
$ws.Range("R9").Value = 70.74722558771001
$ws.Range("S9").Value = 0.06875025630041828
$ws.Range("T9").Value = 0.06875025630041828
$ws.Range("G10").Value = 5.662136666666666
$ws.Range("H10").Value = 16.98641
$ws.Range("I10").Value = 0.129641427533958
$ws.Range("J10").Value = 0.129641427533958
$ws.Range("M10").Value = 0.8675673333333332
$ws.Range("N10").Value = 2.602702
$ws.Range("O10").Value = 0.331395915889395
$ws.Range("P10").Value = 0.331395915889395
$ws.Range("Q10").Value = 4.912284808868888
$ws.Range("R10").Value = 44.21056327981999
$ws.Range("S10").Value = 0.04296263961482465
$ws.Range("T10").Value = 0.04296263961482465
